# novas aulas credit scoring
# Update the confidence-interval and hypothesis-test worksheets with new
# exercise inputs, then leave the UI selection state matching the author's
# last edits (TH-Média_z ends up the active sheet/tab).

$wb = $excel.ActiveWorkbook

# --- IC-Média_Z (sheet 2): confidence interval for the mean, sigma known ---
$wsICMediaZ = $wb.Worksheets.Item(2)
$wsICMediaZ.Range("C3").Value = 0.95
$wsICMediaZ.Range("C5").Value = 50
$wsICMediaZ.Range("C6").Value = 6
$wsICMediaZ.Activate() | Out-Null
$wsICMediaZ.Range("C7").Select() | Out-Null

# --- IC-Média_t (sheet 3): confidence interval for the mean, sigma unknown (t) ---
$wsICMediaT = $wb.Worksheets.Item(3)
$wsICMediaT.Range("C3").Value = 0.9
$wsICMediaT.Range("C5").Value = 1976
$wsICMediaT.Range("C6").Value = 11
$wsICMediaT.Range("C7").Value = 28
$wsICMediaT.Activate() | Out-Null
$wsICMediaT.Range("C8").Select() | Out-Null

# --- TH-Média_z (sheet 6): hypothesis test for the mean, sigma known ---
$wsTHMediaZ = $wb.Worksheets.Item(6)
$wsTHMediaZ.Range("O3").Value = 0.9
$wsTHMediaZ.Range("R3").Value = 110
$wsTHMediaZ.Range("O7").Value = 20
$wsTHMediaZ.Range("O8").Value = 120

# --- TH-Média_t (sheet 7): just a selection/cursor change ---
$wsTHMediaT = $wb.Worksheets.Item(7)
$wsTHMediaT.Activate() | Out-Null
$wsTHMediaT.Range("R3").Select() | Out-Null

# TH-Média_z is the last sheet the author interacted with, so it ends up the
# active/tabSelected sheet with R4 selected (this also moves tabSelected off
# TH-Prop_Z, sheet 8).
$wsTHMediaZ.Activate() | Out-Null
$wsTHMediaZ.Range("R4").Select() | Out-Null
